# Error Calculations and Plots
# Re-derive the "missing data" mask for the imputation dataset:
#  - Drop the "RM 232" and "SC 92" rows entirely (rows shift up).
#  - Update which cells are treated as missing (blank) vs populated
#    in columns B/D/F for several remaining rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two rows that no longer appear in the data set.
# "RM 232" is row 26; after it is removed, "SC 92" (originally row 28)
# becomes row 27, so delete it next.
$ws.Rows.Item(26).Delete()
$ws.Rows.Item(27).Delete()

# Newly-populated cells (previously marked missing).
$ws.Range("D2").Value = -13.5
$ws.Range("F4").Value = 17.97
$ws.Range("D12").Value = -14.1
$ws.Range("D20").Value = -14
$ws.Range("D21").Value = -14.3
$ws.Range("F23").Value = 16.48
$ws.Range("F29").Value = 18.06
$ws.Range("B30").Value = -19.7
$ws.Range("D31").Value = -13.7
$ws.Range("D33").Value = -14.1

# Cells newly marked as missing (cleared).
$ws.Range("F3").ClearContents()
$ws.Range("F5").ClearContents()
$ws.Range("D6").ClearContents()
$ws.Range("F8").ClearContents()
$ws.Range("D14").ClearContents()
$ws.Range("D22").ClearContents()
$ws.Range("D23").ClearContents()
$ws.Range("F27").ClearContents()
$ws.Range("B32").ClearContents()
